$wb = $excel.ActiveWorkbook

# Mapping of sheet name -> list of (row, newValue) for column F (6)

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 2111
$ws.Cells.Item(7, 6).Value = 7944
$ws.Cells.Item(9, 6).Value = 55
$ws.Cells.Item(13, 6).Value = 1541
$ws.Cells.Item(14, 6).Value = 1314
$ws.Cells.Item(15, 6).Value = 183
$ws.Cells.Item(16, 6).Value = 4021
$ws.Cells.Item(17, 6).Value = 698
$ws.Cells.Item(18, 6).Value = 29
$ws.Cells.Item(19, 6).Value = 1104
$ws.Cells.Item(20, 6).Value = 1231
$ws.Cells.Item(21, 6).Value = 428
$ws.Cells.Item(22, 6).Value = 6219
$ws.Cells.Item(25, 6).Value = 4221
$ws.Cells.Item(27, 6).Value = 1946
$ws.Cells.Item(28, 6).Value = 1167
$ws.Cells.Item(29, 6).Value = 302
$ws.Cells.Item(30, 6).Value = 1032
$ws.Cells.Item(33, 6).Value = 203
$ws.Cells.Item(34, 6).Value = 48
$ws.Cells.Item(37, 6).Value = 503
$ws.Cells.Item(38, 6).Value = 1872
$ws.Cells.Item(39, 6).Value = 108
$ws.Cells.Item(40, 6).Value = 408
$ws.Cells.Item(41, 6).Value = 152
$ws.Cells.Item(42, 6).Value = 1141
$ws.Cells.Item(44, 6).Value = 64
$ws.Cells.Item(45, 6).Value = 37
$ws.Cells.Item(48, 6).Value = 172

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(9, 6).Value = 16
$ws.Cells.Item(12, 6).Value = 371
$ws.Cells.Item(20, 6).Value = 173
$ws.Cells.Item(22, 6).Value = 86
$ws.Cells.Item(30, 6).Value = 271

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 454
$ws.Cells.Item(6, 6).Value = 1559
$ws.Cells.Item(7, 6).Value = 470
$ws.Cells.Item(9, 6).Value = 932
$ws.Cells.Item(10, 6).Value = 1073
$ws.Cells.Item(11, 6).Value = 1264
$ws.Cells.Item(12, 6).Value = 1577

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 454
$ws.Cells.Item(3, 6).Value = 1559
$ws.Cells.Item(6, 6).Value = 470
$ws.Cells.Item(8, 6).Value = 2111
$ws.Cells.Item(9, 6).Value = 7944
$ws.Cells.Item(10, 6).Value = 55
$ws.Cells.Item(11, 6).Value = 932
$ws.Cells.Item(13, 6).Value = 16
$ws.Cells.Item(16, 6).Value = 1541
$ws.Cells.Item(17, 6).Value = 1264
$ws.Cells.Item(18, 6).Value = 1314
$ws.Cells.Item(20, 6).Value = 183
$ws.Cells.Item(21, 6).Value = 1577
$ws.Cells.Item(22, 6).Value = 4021
$ws.Cells.Item(23, 6).Value = 371
$ws.Cells.Item(25, 6).Value = 698
$ws.Cells.Item(26, 6).Value = 29
$ws.Cells.Item(27, 6).Value = 1104
$ws.Cells.Item(28, 6).Value = 1231
$ws.Cells.Item(29, 6).Value = 428
$ws.Cells.Item(30, 6).Value = 6219
$ws.Cells.Item(33, 6).Value = 1946
$ws.Cells.Item(34, 6).Value = 1167
$ws.Cells.Item(35, 6).Value = 302
$ws.Cells.Item(37, 6).Value = 173
$ws.Cells.Item(38, 6).Value = 203
$ws.Cells.Item(39, 6).Value = 86
$ws.Cells.Item(40, 6).Value = 503
$ws.Cells.Item(41, 6).Value = 1872
$ws.Cells.Item(42, 6).Value = 108
$ws.Cells.Item(43, 6).Value = 408
$ws.Cells.Item(44, 6).Value = 1141
$ws.Cells.Item(47, 6).Value = 271
$ws.Cells.Item(49, 6).Value = 172
